$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "选择3个问题" -> "选择3-5个问题", with "3" and the new "-5" run bold
#    and colored C00000 (dark red) instead of the plain 333333 gray.
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("选择3个问题", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$digitStart = $r.Start + 2
$p3 = $d.Range($digitStart, $digitStart + 1)
$p3.InsertAfter("-5")

$p3a = $d.Range($digitStart, $digitStart + 1)
$p3a.Font.Bold = 1
$p3a.Font.Color = 192

$p3b = $d.Range($digitStart + 1, $digitStart + 3)
$p3b.Font.Bold = 1
$p3b.Font.Color = 192

# ---------------------------------------------------------------------
# 2) Merge the "用Git" run and the "进行版本控制" run that follows the
#    (currently misplaced) _GoBack bookmark into a single run, without
#    disturbing the separate "。" run that comes right after.
# ---------------------------------------------------------------------
$g = $d.Content
$null = $g.Find.Execute("用Git", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gitStart = $g.Start
$gitEnd = $g.End
$gitRange = $d.Range($gitStart, $gitEnd)
$gitRange.Text = "用Git进行版本控制"

$afterLen = $gitStart + [string]"用Git进行版本控制".Length
$stale = $d.Content
$stale.Start = $afterLen
$null = $stale.Find.Execute("进行版本控制", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$stale.Delete()

# ---------------------------------------------------------------------
# 3) Relocate the (hidden) _GoBack bookmark to sit right after the final
#    "。" of the second paragraph ("综合实践：…社区编程等。"), which both
#    removes it from between the runs merged above and places it at its
#    new home. Bookmarks.Add with an existing name moves the bookmark.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$marker = $d.Range($p2.End - 1, $p2.End - 1)
$marker.InsertAfter("_TmpBookmarkAnchor_")

$anchor = $d.Content
$null = $anchor.Find.Execute("_TmpBookmarkAnchor_", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $anchor)
$anchor.Text = ""
